$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attendance")

# Mark attendance ("P" = present) for the first lab (column C) for the
# students belonging to the second semi-group.
$ws.Range("C2").Value = "P"

# Rename a few students to indicate the secondary group ("gr 4") they
# actually attend, per the commit "Lab 01: prezenta semigrupa 2."
$ws.Range("B24").Value = "Vranau V. Flavius Silviu (gr 4)"
$ws.Range("B6").Value  = "Dulau I. Marius Cristian (joi gr 4 sem 1)"
$ws.Range("B18").Value = "Prata L. Dragos Liviu (gr 4)"

$presentRows = 6,14,16,17,18,19,20,22,23,24
foreach ($r in $presentRows) {
    $ws.Range("C$r").Value = "P"
}

$ws.Range("Q18").Select()
